$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Trends Status" (sheet1): update counts / percentages
# ---------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("B2").Value = 0
$wsTrends.Range("C2").Value = 4
$wsTrends.Range("D2").Value = 0
$wsTrends.Range("E2").Value = 8.300000000000001

$wsTrends.Range("C3").Value = 12
$wsTrends.Range("E3").Value = 25

$wsTrends.Range("B4").Value = 1
$wsTrends.Range("C4").Value = 26
$wsTrends.Range("D4").Value = 100
$wsTrends.Range("E4").Value = 54.2

$wsTrends.Range("C5").Value = 2
$wsTrends.Range("E5").Value = 4.2

$wsTrends.Range("C6").Value = 4
$wsTrends.Range("E6").Value = 8.300000000000001

$wsTrends.Range("B7").Value = 13
$wsTrends.Range("C7").Value = 56

$wsTrends.Range("B8").Value = 379
$wsTrends.Range("C8").Value = 289

# ---------------------------------------------------------------------
# Sheet "Priority Status" (sheet3): update species counts
# ---------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------
# Sheet "Species qualification" (sheet4): update label + counts
# ---------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B2").Value = 393

$wsQual.Range("B3").Value = 14
$wsQual.Range("C3").Value = 1

$wsQual.Range("B4").Value = 104
$wsQual.Range("C4").Value = 48

# ---------------------------------------------------------------------
# Duplicate "High Priority break-up" (sheet5) into a new sheet placed
# right after it, BEFORE changing any values, so the new sheet keeps
# the original (pre-update) figures.
# ---------------------------------------------------------------------
$wsHighPriority = $wb.Worksheets.Item("High Priority break-up")

$wsMajorUpdate = $wb.Worksheets.Add($null, $wsHighPriority)
$wsMajorUpdate.Name = "Major update - High Priority "

$wsMajorUpdate.Range("A1").Value = "Break-up"
$wsMajorUpdate.Range("B1").Value = "High Species (no.)"
$wsMajorUpdate.Range("C1").Value = "High Species (perc.)"
$wsMajorUpdate.Range("D1").Value = "New High Species (no.)"
$wsMajorUpdate.Range("E1").Value = "New High Species (perc.)"
$wsMajorUpdate.Range("A1:E1").Font.Bold = $true
$wsMajorUpdate.Range("A1:E1").HorizontalAlignment = -4108

$wsMajorUpdate.Range("A2").Value = "Trend New"
$wsMajorUpdate.Range("B2").Value = 1
$wsMajorUpdate.Range("C2").Value = 8.300000000000001
$wsMajorUpdate.Range("D2").Value = 1
$wsMajorUpdate.Range("E2").Value = 8.300000000000001

$wsMajorUpdate.Range("A3").Value = "IUCN"
$wsMajorUpdate.Range("B3").Value = 11
$wsMajorUpdate.Range("C3").Value = 91.7
$wsMajorUpdate.Range("D3").Value = 11
$wsMajorUpdate.Range("E3").Value = 91.7

# ---------------------------------------------------------------------
# Now update the original "High Priority break-up" sheet's figures and
# rename it to "Interannual update - High Pri".
# ---------------------------------------------------------------------
$wsHighPriority.Range("B2").Value = 72
$wsHighPriority.Range("C2").Value = 69.90000000000001
$wsHighPriority.Range("D2").Value = 72
$wsHighPriority.Range("E2").Value = 77.40000000000001

$wsHighPriority.Range("B3").Value = 31
$wsHighPriority.Range("C3").Value = 30.1
$wsHighPriority.Range("D3").Value = 21
$wsHighPriority.Range("E3").Value = 22.6

$wsHighPriority.Name = "Interannual update - High Pri"
